# Apply text corrections (typo/spelling fixes, trailing punctuation/whitespace
# cleanups) to the "Inital_CETSA_clusters" worksheet, and update the saved
# window selection/scroll position, matching the upstream "Add files via
# upload" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B ("Function") fixes
$ws.Range("B2").Value  = "Key protein in glutathione synthesis"
$ws.Range("B15").Value = "Move ubiqutin from E1 to E3 ligases"
$ws.Range("B18").Value = "Aifm1 is release upon apoptosis to enter nucleus where it acivates DNAses initating DNA fragmentation"
$ws.Range("B21").Value = "Chek1 signals ss-DNA damage, DNMT1 transfer DNA metylation "
$ws.Range("B24").Value = "Transcription factors, often in immune cell contexts"

# Column A ("Name of cluster") fixes
$ws.Range("A7").Value  = "Cu/Zn metabolism"
$ws.Range("A11").Value = "Trifunctional enzyme complex - lipid beta oxidation"
$ws.Range("A19").Value = "PI3K pathway induced autophagy"
$ws.Range("A23").Value = "Translesion synthesis (TLS) signature"
$ws.Range("A25").Value = "TNFa- NF-kB signalling"
$ws.Range("A29").Value = "late G2, M entry  - change in cell cycle distribution."
$ws.Range("A31").Value = "Sugar matabolism downstream of PI3K"
$ws.Range("A36").Value = "Chek2 signaling - double strand breaks"

# Column D ("Functional hypothesis for shifts") fixes
$ws.Range("D2").Value  = "Flux and allosteric regulation. Stability shifts of regulatory GCLM subunit might report on activation state."
$ws.Range("D4").Value  = "Mainly nuclear matrix localised, can be induce by effector casapse attack"
$ws.Range("D8").Value  = "Decrease on multiple proteins reflect cellullar thymidine level decrease"
$ws.Range("D10").Value = "Destabilizations can reflect release from DNA, eg to be replaced by repair polymerases to "
$ws.Range("D12").Value = " Stabilizations and/or level increase reflect CDK activation. Stability shift both increase interactions and phosorylations"
$ws.Range("D20").Value = "Stabilization likely reflects binding to DNA"
$ws.Range("D24").Value = "Could be phosphorylation and DEAR effect when entering nucleus"
$ws.Range("D25").Value = "NFKBIB  degraded upon pathway activation. Stabilization of NFKB1 due to activation by phosphorylation"

# Update the saved view: scroll so row 19 is at the top and select A36,
# matching the workbook's recorded sheetView/selection state.
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A36").Select()
